# Bug fix in Eduati data files (DIFI_noCTRL_meas.xlsx)
#
# Sheet1 ("DIFI_noCTRL_meas") had 43 stray leftover rows (45-87) that only
# contained a running index in column A (left over from a larger source
# sheet). Those rows are removed so the sheet's real data block (A1:N44)
# is the whole sheet. The previously-active tab/selection state also gets
# reset: Sheet1 (rather than Sheet3) is now the active/selected tab, with
# the cursor on I49, and Sheet3's own selection is left as-is (A2:N44)
# while it is no longer the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet1: drop the 43 stray rows below the real data (rows 45-87) ---
# This shrinks the sheet's dimension from A1:N87 back down to A1:N44.
$ws1.Rows("45:87").Delete()

# --- Restore Sheet3's own (unchanged) selection before we move away from it ---
$ws3.Activate()
$ws3.Range("A2:N44").Select()
# Sheet3 is scrolled down a bit (top-left visible row 19) even though the
# selection itself stays A2:N44.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet1 becomes the active/selected worksheet, cursor parked at I49 ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("I49").Select()
